$wb = $excel.ActiveWorkbook

function Add-Row89($ws, $bVal, $cVal, $dVal, $eVal, $fVal, $gVal, $hVal, $iVal) {
    $r = 89
    $ws.Cells.Item($r, 1).Value = 45875.49188657408
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $bVal
    $ws.Cells.Item($r, 3).Value = $cVal
    $ws.Cells.Item($r, 4).Value = $dVal
    $ws.Cells.Item($r, 5).Value = $eVal
    $ws.Cells.Item($r, 6).Value = $fVal
    $ws.Cells.Item($r, 7).Value = [double]$gVal
    $ws.Cells.Item($r, 8).Value = $hVal
    $ws.Cells.Item($r, 9).Value = $iVal
}

# FE_LFT_#1
$ws1 = $wb.Worksheets.Item(1)
Add-Row89 $ws1 "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x01,0x10" "0xf" 380 "7.598631275147109e+23" 272 15

# FE_LFT_#2
$ws2 = $wb.Worksheets.Item(2)
Add-Row89 $ws2 "0x01,0x90" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x01,0x20" "0xe" 400 "5.68432987514711e+23" 288 14

# FE_PLT_#1
$ws3 = $wb.Worksheets.Item(3)
Add-Row89 $ws3 "0x00,0x6e" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x5D" "0x3" 110 "5.68631262647114e+23" 93 3

# FE_PLT_#2
$ws4 = $wb.Worksheets.Item(4)
Add-Row89 $ws4 "0x00,0x6e" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x5C" "0x3" 110 "9.85046333984776e+23" 92 3
